# Apply "More report details, UI update" changes to commData sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G (MaxDistance): change every data row (2-30) from 1000 -> 12
for ($r = 2; $r -le 30; $r++) {
    $ws.Cells.Item($r, 7).Value = 12
}

# Column F (AffectedPop): rows 2 and 3 change from 10 -> 0
$ws.Cells.Item(2, 6).Value = 0
$ws.Cells.Item(3, 6).Value = 0

# Column A (Active, boolean): toggle specific rows
$ws.Cells.Item(2, 1).Value = $true
$ws.Cells.Item(6, 1).Value = $false
$ws.Cells.Item(10, 1).Value = $true
$ws.Cells.Item(11, 1).Value = $true
$ws.Cells.Item(14, 1).Value = $true
$ws.Cells.Item(21, 1).Value = $false
